# Commit message:
#   "Updated the script so that the bacteria subset keep their bioproj
#    reference. No need to name them 'Otu. the .txt and .xlsx is updated"
#
# This script:
#   1. Renames the query/worksheet "my_file" -> "mock_genomes"
#      (the workbook's defined name ExternalData_1 follows automatically).
#   2. Renames the query table ("my_file") to "mock_genomes" and switches
#      its visual style from TableStyleLight11 to TableStyleLight1.
#   3. Replaces the bacteria/archaea genome_id values (rows 37-71 of
#      column A) from the placeholder "Otu###" labels with their real
#      BioProject accession numbers ("PRJNA######"), in row order.
#   4. Moves the active cell selection on the data sheet to H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("my_file")   # the "my_file" query sheet (currently the active tab)

# --- 1. Rename the worksheet (defined names referencing it update too) ---
$ws.Name = "mock_genomes"

# --- 2. Rename the query table and change its style ---
$lo = $ws.ListObjects.Item(1)
$lo.Name = "mock_genomes"
$lo.TableStyle = "TableStyleLight1"

# --- 3. Replace the bacteria/archaea genome_id values with BioProject IDs ---
$genomeIds = @(
    "PRJNA67115",
    "PRJNA217481",
    "PRJNA20011",
    "PRJNA186462",
    "PRJNA33691",
    "PRJNA33599",
    "PRJNA20399",
    "PRJNA212980",
    "PRJNA261945",
    "PRJNA183309",
    "PRJNA27951",
    "PRJNA261104",
    "PRJNA183018",
    "PRJNA238302",
    "PRJNA42009",
    "PRJNA80827",
    "PRJNA171367",
    "PRJNA232079",
    "PRJNA242829",
    "PRJNA165395",
    "PRJNA81617",
    "PRJNA13473",
    "PRJNA186910",
    "PRJNA63851",
    "PRJNA18505",
    "PRJNA20391",
    "PRJNA12634",
    "PRJNA17707",
    "PRJNA232351",
    "PRJNA42475",
    "PRJNA168616",
    "PRJNA182711",
    "PRJNA15771",
    "PRJNA256039",
    "PRJNA190819"
)

$startRow = 37
for ($i = 0; $i -lt $genomeIds.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $genomeIds[$i]
}

# --- 4. Update the active selection on the data sheet ---
$ws.Activate()
$ws.Range("H7").Select()
